$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text edits (edit only the changed run, preserve the rest of the rich text) ----
$ws.Range("A8").Characters(21, 1).Text = "8"
$ws.Range("C9").Characters(27, 9).Text = "2/19/2024"
$ws.Range("C9").Characters(47, 9).Text = "2/25/2024"

# ---- Template cells, used via Copy() to clone an exact cell style (avoids NumberFormat
# assignment, which always mints a brand-new style id in this runtime) ----
$tmplText0  = $ws.Range("C14")   # s=14 t=s -> shared string idx 20 ("0")
$tmplTextNA = $ws.Range("E14")   # s=14 t=s -> shared string idx 21 ("***.*")
$tmplNum15  = $ws.Range("I15")   # s=15 numeric template
$tmplNum16  = $ws.Range("K15")   # s=16 numeric template

# ---- Row 15 ----

# ---- Row 15 ----
$tmplText0.Copy($ws.Range("C15"))
$tmplText0.Copy($ws.Range("D15"))
$tmplTextNA.Copy($ws.Range("E15"))
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("L15").Value = -66.666666666666

# ---- Row 16 ----
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -40.90909090909
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 38
$ws.Range("K16").Value = -31.578947368421
$ws.Range("L16").Value = 36.842105263157
$ws.Range("M16").Value = -25.714285714285
$ws.Range("N16").Value = -85.393258426966

# ---- Row 17 ----
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -19.047619047619
$ws.Range("I17").Value = 38
$ws.Range("J17").Value = 39
$ws.Range("K17").Value = -2.564102564102
$ws.Range("L17").Value = 35.714285714285
$ws.Range("M17").Value = 80.95238095238
$ws.Range("N17").Value = -50

# ---- Row 18 ----
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 13.333333333333
$ws.Range("I18").Value = 61
$ws.Range("J18").Value = 39
$ws.Range("K18").Value = 56.410256410256
$ws.Range("L18").Value = 24.489795918367
$ws.Range("M18").Value = -10.294117647058
$ws.Range("N18").Value = -63.905325443787

# ---- Row 19 ----
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 75
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = -12.698412698412
$ws.Range("I19").Value = 97
$ws.Range("J19").Value = 130
$ws.Range("K19").Value = -25.384615384615
$ws.Range("L19").Value = 12.790697674418
$ws.Range("M19").Value = 61.666666666666
$ws.Range("N19").Value = 59.016393442622

# ---- Row 20 ----
$tmplNum15.Copy($ws.Range("C20"))
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 19
$ws.Range("K20").Value = -10.526315789473
$ws.Range("L20").Value = -26.086956521739
$ws.Range("M20").Value = -22.727272727272
$ws.Range("N20").Value = -86.71875

# ---- Row 21 ----
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 133
$ws.Range("H21").Value = -15.78947368421
$ws.Range("I21").Value = 241
$ws.Range("J21").Value = 271
$ws.Range("K21").Value = -11.070110701107
$ws.Range("L21").Value = 15.865384615384
$ws.Range("M21").Value = 15.865384615384
$ws.Range("N21").Value = -61.254019292604

# ---- Row 22 ----
$tmplNum15.Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$tmplNum16.Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100
$tmplText0.Copy($ws.Range("F22"))
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = -50

# ---- Row 23 ----
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 133.333333333333
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 25
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = -3.846153846153
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 47.058823529411

# ---- Row 24 ----
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 81.25
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 66
$ws.Range("H24").Value = 45.454545454545
$ws.Range("I24").Value = 181
$ws.Range("J24").Value = 155
$ws.Range("K24").Value = 16.774193548387
$ws.Range("L24").Value = 10.365853658536
$ws.Range("M24").Value = 7.100591715976

# ---- Row 25 ----
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 87.5
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = -5.555555555555
$ws.Range("I25").Value = 79
$ws.Range("J25").Value = 68
$ws.Range("K25").Value = 16.176470588235
$ws.Range("L25").Value = 25.396825396825
$ws.Range("M25").Value = 27.419354838709

# ---- Row 26 ----
$tmplText0.Copy($ws.Range("D26"))
$tmplTextNA.Copy($ws.Range("E26"))
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -50
$ws.Range("L26").Value = -60

# ---- Row 27 ----
$tmplNum15.Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = 100

# ---- Row 30 ----
$ws.Range("L30").Value = -37.5
